# Applies the edits described by the diff:
#  - workbook window geometry + activeTab change
#  - new shared string "f_agg_Category"
#  - sheet "_set_FLOWS": drop tab-selected, new selection, set column widths
#  - sheet "_set_FLOWS_AGG": becomes the active/selected tab, new selection,
#    set column widths, rename header B1 from "f_Category" to "f_agg_Category"

$wb = $excel.ActiveWorkbook

$wsFlows = $wb.Worksheets.Item("_set_FLOWS")
$wsAgg   = $wb.Worksheets.Item("_set_FLOWS_AGG")

# --- Rename header on _set_FLOWS_AGG: B1 "f_Category" -> "f_agg_Category" ---
$wsAgg.Range("B1").Value = "f_agg_Category"

# --- Column widths (bestFit / autofit-style custom widths) ---
$wsFlows.Columns.Item(1).ColumnWidth = 12.7265625
$wsFlows.Columns.Item(2).ColumnWidth = 16
$wsFlows.Columns.Item(3).ColumnWidth = 12.453125

$wsAgg.Columns.Item(1).ColumnWidth = 11.1796875
$wsAgg.Columns.Item(2).ColumnWidth = 16

# --- Selections on each sheet ---
$wsFlows.Range("C31").Select()
$wsAgg.Range("U8").Select()

# --- Make _set_FLOWS_AGG the active/selected tab ---
$wsAgg.Activate()
$wsAgg.Select()
